# ----------------------------------------------------------------------------
# Rework the "Yearly" sheet header block: add a Market Cap column, split the
# P/E-family headers (PEG / Current P/E(x) / EPS / PE), add an Expenses line to
# every yearly block, relabel several headers with their unit suffixes, and
# resize the grid accordingly (mirrors the switch to scraping with
# mechanicalsoup, which returns a wider set of columns per filing).
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Yearly")
$ws2 = $wb.Worksheets.Item("Quarterly")

# Full header row (row 5), column A (1) through CZ (104).
$headerTexts = @(
  "Company Name", "Sector/Industry", "Group", "Market Cap", "CMP", "Valuation", "PEG", "Current P/E(x)",
  "Avg PE (5 yrs)", "EPS", "EPS CAGR %(3 yrs)", "ROE(%)", "ROCE(%)", "P/B(x)", "Dividend Yield", "Debt to equity ratio",
  "Interest Coverage Ratios", "Revenue CAGR(3 yrs)", "Profit CAGR(3 yrs)", "Net Profit CAGR(3 yrs)", "Revenue", "Expenses", "Profit", "Net Profit",
  "Net Cash Flow", "EPS", "PE", "Revenue", "Expenses", "Profit", "Net Profit", "Net Cash Flow",
  "EPS", "PE", "Revenue", "Expenses", "Profit", "Net Profit", "Net Cash Flow", "EPS",
  "PE", "Revenue", "Expenses", "Profit", "Net Profit", "Net Cash Flow", "EPS", "PE",
  "Revenue", "Expenses", "Profit", "Net Profit", "Net Cash Flow", "EPS", "PE", "Revenue",
  "Expenses", "Profit", "Net Profit", "Net Cash Flow", "EPS", "PE", "Revenue", "Expenses",
  "Profit", "Net Profit", "Net Cash Flow", "EPS", "PE", "Revenue", "Expenses", "Profit",
  "Net Profit", "Net Cash Flow", "EPS", "PE", "Revenue", "Expenses", "Profit", "Net Profit",
  "Net Cash Flow", "EPS", "PE", "Revenue", "Expenses", "Profit", "Net Profit", "Net Cash Flow",
  "EPS", "PE", "Revenue", "Expenses", "Profit", "Net Profit", "Net Cash Flow", "EPS",
  "PE", "Revenue", "Expenses", "Profit", "Net Profit", "Net Cash Flow", "EPS", "PE"
)
for ($i = 0; $i -lt $headerTexts.Length; $i++) {
  $ws.Cells.Item(5, $i + 1).Value = $headerTexts[$i]
}

# The last two columns of the block (EPS / PE for the oldest period) close the
# table visually, so they reuse the bordered "closing" style that previously
# sat on the final column of the old, narrower layout.
$ws.Range("BS5").Copy()
$ws.Range("CY5:CZ5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 5 now wraps a bit more than the sheet default, matching the height Excel
# computes automatically for the refreshed header.
$ws.Rows.Item(5).RowHeight = 15

# Column widths, recalculated (AutoFit) for the new header text / column count.
$colWidths = @(
  25.85546875, 15.0, 6.7109375, 11.28515625, 5.28515625, 9.7109375, 4.7109375, 13.85546875,
  12.85546875, 4.28515625, 17.42578125, 7.7109375, 8.85546875, 6.7109375, 14.28515625, 18.7109375,
  23.140625, 20.140625, 17.140625, 21.0, 9.0, 9.42578125, 11.42578125, 9.85546875,
  13.85546875, 4.28515625, 3.28515625, 9.0, 9.42578125, 6.140625, 9.85546875, 13.85546875,
  4.28515625, 3.28515625, 9.0, 9.42578125, 6.140625, 9.85546875, 13.85546875, 4.28515625,
  3.28515625, 9.0, 9.42578125, 6.140625, 9.85546875, 13.85546875, 4.28515625, 3.28515625,
  9.0, 9.42578125, 6.140625, 9.85546875, 13.85546875, 4.28515625, 3.28515625, 9.0,
  9.42578125, 6.140625, 9.85546875, 13.85546875, 4.28515625, 3.28515625, 9.0, 9.42578125,
  6.140625, 9.85546875, 13.85546875, 4.28515625, 3.28515625, 9.0, 9.42578125, 6.140625,
  9.85546875, 13.85546875, 4.28515625, 3.28515625, 9.0, 9.42578125, 6.140625, 9.85546875,
  13.85546875, 4.28515625, 3.28515625, 9.0, 9.42578125, 6.140625, 9.85546875, 13.85546875,
  4.28515625, 3.28515625, 9.0, 9.42578125, 6.140625, 9.85546875, 13.85546875, 4.28515625,
  3.28515625, 9.0, 9.42578125, 6.140625, 9.85546875, 13.85546875, 4.28515625, 3.28515625
)
for ($i = 0; $i -lt $colWidths.Length; $i++) {
  $ws.Columns.Item($i + 1).ColumnWidth = $colWidths[$i]
}

# Selection reverts to the whole-sheet default (no more sticky A6:BV7 filter
# selection left over from the old scraper).
$ws.Cells.Select()

# "Quarterly" sheet note references the same shared string, just keep it in sync.
$ws2.Range("A2").Value = "Same structure as yearly"

$wb.Save()
